$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Update SKU ID (col B) and Quantity (col C) values for rows 2-14
$ws.Range("B2").Value = "11-GTZ-7VMF"
$ws.Range("C2").Value = 7

$ws.Range("B3").Value = "11-H12V3530H5"
$ws.Range("C3").Value = 7

$ws.Range("B4").Value = "12-100020-00000"
$ws.Range("C4").Value = 9

$ws.Range("B5").Value = "12-100020-FLAP"
$ws.Range("C5").Value = 5

$ws.Range("B6").Value = "12-50012-00000"
$ws.Range("C6").Value = 27

$ws.Range("B7").Value = "FP-06435-KPP-2700"
$ws.Range("C7").Value = 22

$ws.Range("B8").Value = "FP-43120-362-2700"
$ws.Range("C8").Value = 140

$ws.Range("B9").Value = "FP-43125-KGA-2700"
$ws.Range("C9").Value = 140

$ws.Range("B10").Value = "FP-54410-THU-2700"
$ws.Range("C10").Value = 140

$ws.Range("B11").Value = "FP-F533A-RXK-2700"
$ws.Range("C11").Value = 140

$ws.Range("B12").Value = "FP-W0045-F1Z-2700"
$ws.Range("C12").Value = 160

$ws.Range("B13").Value = "GSMF-GTZ-4V"
$ws.Range("C13").Value = 11

$ws.Range("B14").Value = "H2-231PA-KZL-1200"
$ws.Range("C14").Value = 1

# Remove the former row 15 (data now ends at row 14)
$ws.Rows.Item(15).Delete()
